# This script applies corrected IFRS-list figures onto rows 2-9 (company_list sheet)
# following the upstream "error solve ifrs list" commit: updates every numeric
# metric cell to the corrected value and clears cells that no longer apply.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values, keyed by A1 reference
$newValues = @{
    "D2" = 2878
    "E2" = 117
    "F2" = 117
    "G2" = 110
    "H2" = 88
    "I2" = 88
    "K2" = 1032
    "L2" = 640
    "M2" = 393
    "N2" = 393
    "P2" = 62
    "Q2" = 29
    "R2" = -51
    "S2" = 24
    "T2" = 45
    "U2" = -16
    "V2" = 191
    "W2" = 4.06
    "X2" = 3.06
    "Y2" = 25.85
    "Z2" = 9.31
    "AA2" = 163.01
    "AB2" = 495.37
    "AC2" = 708
    "AD2" = 26.9
    "AE2" = 3172
    "AF2" = 6.01
    "AG2" = 75
    "AH2" = 0.39
    "AI2" = 10.59
    "AJ2" = 11558200
    "D3" = 3158
    "E3" = 111
    "F3" = 111
    "G3" = 94
    "H3" = 74
    "I3" = 74
    "K3" = 1202
    "L3" = 696
    "M3" = 506
    "N3" = 506
    "P3" = 62
    "Q3" = 96
    "R3" = -107
    "S3" = 11
    "T3" = 63
    "U3" = 33
    "V3" = 220
    "W3" = 3.53
    "X3" = 2.35
    "Y3" = 16.51
    "Z3" = 6.64
    "AA3" = 137.44
    "AB3" = 575.88
    "AC3" = 597
    "AD3" = 24.97
    "AE3" = 4106
    "AF3" = 3.63
    "AG3" = 100
    "AH3" = 0.67
    "AI3" = 16.68
    "AJ3" = 11558200
    "D4" = 3404
    "E4" = 113
    "F4" = 113
    "G4" = 100
    "H4" = 73
    "I4" = 73
    "K4" = 1224
    "L4" = 692
    "M4" = 531
    "N4" = 531
    "P4" = 62
    "Q4" = 128
    "R4" = -7
    "S4" = -56
    "T4" = 66
    "U4" = 63
    "V4" = 213
    "W4" = 3.32
    "X4" = 2.15
    "Y4" = 14.12
    "Z4" = 6.04
    "AA4" = 130.35
    "AB4" = 648.95
    "AC4" = 589
    "AD4" = 19.85
    "AE4" = 4386
    "AF4" = 2.67
    "AG4" = 200
    "AH4" = 1.71
    "AI4" = 33.13
    "AJ4" = 11558200
    "D5" = 3950
    "E5" = 122
    "F5" = 122
    "G5" = 118
    "H5" = 95
    "I5" = 95
    "K5" = 1475
    "L5" = 851
    "M5" = 624
    "N5" = 624
    "P5" = 62
    "Q5" = 45
    "R5" = -114
    "S5" = 7
    "T5" = 74
    "U5" = -30
    "V5" = 234
    "W5" = 3.09
    "X5" = 2.41
    "Y5" = 16.48
    "Z5" = 7.06
    "AA5" = 136.22
    "AB5" = 740.17
    "AC5" = 766
    "AD5" = 16.44
    "AE5" = 5150
    "AF5" = 2.45
    "AG5" = 500
    "AH5" = 3.97
    "AI5" = 63.71
    "AJ5" = 11558200
    "D6" = 4021
    "E6" = 113
    "F6" = 113
    "G6" = 91
    "H6" = 75
    "I6" = 75
    "K6" = 1214
    "L6" = 736
    "M6" = 478
    "N6" = 478
    "P6" = 62
    "Q6" = 204
    "R6" = -57
    "S6" = -142
    "T6" = 60
    "U6" = 143
    "V6" = 161
    "W6" = 2.82
    "X6" = 1.87
    "Y6" = 13.65
    "Z6" = 5.6
    "AA6" = 154.15
    "AB6" = 721.95
    "AC6" = 605
    "AD6" = 17.18
    "AE6" = 3959
    "AF6" = 2.63
    "AG6" = 500
    "AH6" = 4.81
    "AI6" = 80.25
    "AJ6" = 11558200
    "D7" = 3935
    "E7" = 130
    "G7" = 110
    "H7" = 90
    "I7" = 90
    "W7" = 3.3
    "X7" = 2.29
    "AC7" = 724
    "AD7" = 21.27
    "D8" = 4208
    "E8" = 160
    "G8" = 150
    "H8" = 120
    "I8" = 130
    "W8" = 3.8
    "X8" = 2.85
    "AC8" = 1046
    "AD8" = 14.72
}
foreach ($ref in $newValues.Keys) {
    $ws.Range($ref).Value = $newValues[$ref]
}

# Cells that are no longer populated for these rows (cleared, not zeroed)
$clearedRefs = @(
    "J2",
    "O2",
    "J3",
    "O3",
    "J4",
    "O4",
    "J5",
    "O5",
    "K7",
    "L7",
    "M7",
    "N7",
    "P7",
    "Q7",
    "R7",
    "S7",
    "T7",
    "U7",
    "Y7",
    "Z7",
    "AA7",
    "AE7",
    "AF7",
    "AG7",
    "AH7",
    "AI7",
    "K8",
    "L8",
    "M8",
    "N8",
    "P8",
    "Q8",
    "R8",
    "S8",
    "T8",
    "U8",
    "Y8",
    "Z8",
    "AA8",
    "AE8",
    "AF8",
    "AG8",
    "AH8",
    "AI8",
    "D9",
    "E9",
    "G9",
    "H9",
    "I9",
    "K9",
    "L9",
    "M9",
    "N9",
    "P9",
    "Q9",
    "R9",
    "S9",
    "T9",
    "U9",
    "W9",
    "X9",
    "Y9",
    "Z9",
    "AA9",
    "AC9",
    "AD9",
    "AE9",
    "AF9",
    "AG9",
    "AH9",
    "AI9"
)
foreach ($ref in $clearedRefs) {
    $ws.Range($ref).ClearContents()
}
